# Auto-generated PowerShell COM-interop script
# Applies updated 'want-to-go count' (F) / 'min ticket price' (G) values
# scraped at a later point in time, per commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 286
$ws.Range("F3").Value = 475
$ws.Range("F4").Value = 537
$ws.Range("G4").Value = 78
$ws.Range("F5").Value = 2602
$ws.Range("G5").Value = 90
$ws.Range("F7").Value = 96
$ws.Range("G7").Value = 40
$ws.Range("F8").Value = 102
$ws.Range("G8").Value = 40
$ws.Range("F9").Value = 1733
$ws.Range("G9").Value = 78
$ws.Range("F10").Value = 1733
$ws.Range("G10").Value = 78
$ws.Range("F11").Value = 1438
$ws.Range("F16").Value = 1048
$ws.Range("F19").Value = 261
$ws.Range("F20").Value = 7668
$ws.Range("F21").Value = 8819
$ws.Range("F29").Value = 26
$ws.Range("F31").Value = 367
$ws.Range("F33").Value = 33
$ws.Range("F37").Value = 309
$ws.Range("F38").Value = 43
$ws.Range("F39").Value = 835
$ws.Range("F42").Value = 381
$ws.Range("F44").Value = 234
$ws.Range("F48").Value = 209
$ws.Range("F49").Value = 60

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 2
$ws.Range("F15").Value = 30
$ws.Range("F20").Value = 333

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 2674
$ws.Range("G6").Value = 138

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 286
$ws.Range("F4").Value = 475
$ws.Range("F8").Value = 537
$ws.Range("G8").Value = 78
$ws.Range("F9").Value = 2602
$ws.Range("G9").Value = 90
$ws.Range("F11").Value = 96
$ws.Range("G11").Value = 40
$ws.Range("F12").Value = 102
$ws.Range("G12").Value = 40
$ws.Range("F13").Value = 1733
$ws.Range("G13").Value = 78
$ws.Range("F14").Value = 1733
$ws.Range("G14").Value = 78
$ws.Range("F18").Value = 1048
$ws.Range("F23").Value = 261
$ws.Range("F24").Value = 7668
$ws.Range("F25").Value = 7668
$ws.Range("F26").Value = 8819
$ws.Range("F33").Value = 367
$ws.Range("F35").Value = 33
$ws.Range("F39").Value = 309
$ws.Range("F40").Value = 43
$ws.Range("F41").Value = 835
$ws.Range("F43").Value = 381
$ws.Range("F45").Value = 234
$ws.Range("F49").Value = 209
$ws.Range("F50").Value = 333
$ws.Range("F51").Value = 60
